$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$cell = $table.Cell(3, 2)
$cellRng = $cell.Range
$cellRng.Find.Execute("15:30", $false, $false, $false, $false, $false, $true, 0, $false, "16:00", 2) | Out-Null
Write-Output "Found: $($cellRng.Find.Found)"

$cell2 = $table.Cell(3, 2)
$cellRng2 = $cell2.Range
Write-Output "fresh cell text: [$($cellRng2.Text)]"
Write-Output "fresh cell start/end: $($cellRng2.Start) $($cellRng2.End)"

# Place bookmark right after "16:00" (before the end-of-cell marker)
$bmPos = $cellRng2.Start + 5
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
Write-Output "Bookmark added"
